# Add newly "deleted" SAM variables to the "SAM Variable Changes" sheet.
# These rows document variables removed from Molten Salt Tower Power Block,
# Molten Salt Tower Receiver, Molten Salt Power Block, and Molten Salt Tower
# Storage modules that were redundant, unused, or constant.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SAM Variable Changes")

# Column layout (matches existing rows 2-68):
#   A = Type ("Deleted variable")
#   B = Variable Type ("number")
#   C = Old Name
#   E = Input Page
#   F = Default Value (if new) or reason deleted ("not used")
#   G = Handled in Version Upgrader? ("N/A")
#   H = "Ty"

$rows = @(
  @("m_dot_htf_ref",                  "Molten Salt Tower Power Block"),
  @("T_pb_out",                       "Molten Salt Tower Power Block"),
  @("mode",                           "Molten Salt Tower Power Block"),
  @("fthr_ok",                        "Molten Salt Tower Power Block"),
  @("pb_fixed_par_cntl",              "Molten Salt Tower Power Block"),
  @("dt_cold",                        "Molten Salt Tower Power Block"),
  @("dt_hot",                         "Molten Salt Tower Power Block"),
  @("hx_config",                      "Molten Salt Tower Power Block"),
  @("is_hx",                          "Molten Salt Tower Power Block"),
  @("tech_type",                      "Molten Salt Tower Power Block"),
  @("deg_wind",                       "Molten Salt Tower Receiver"),
  @("P_htf",                          "Molten Salt Tower Receiver"),
  @("T_salt_cold",                    "Molten Salt Power Block"),
  @("HTF",                            "Molten Salt Power Block"),
  @("Design_power",                   "Molten Salt Power Block"),
  @("csp.pt.pwrb.min_restart_time",   "Molten Salt Power Block"),
  @("csp.pt.rec.max_rec_flux",        "Molten Salt Power Block"),
  @("store_fluid",                    "Molten Salt Tower Storage")
)

$startRow = 69
$lastExisting = 68
$r = $startRow
foreach ($row in $rows) {
  $oldName = $row[0]
  $inputPage = $row[1]

  # Clone formatting (styles, fills, borders) from the last existing row.
  $srcRange = $ws.Range("A" + $lastExisting + ":H" + $lastExisting)
  $dstRange = $ws.Range("A" + $r + ":H" + $r)
  $srcRange.Copy($dstRange)

  $ws.Cells.Item($r, 1).Value = "Deleted variable"
  $ws.Cells.Item($r, 2).Value = "number"
  $ws.Cells.Item($r, 3).Value = $oldName
  $ws.Cells.Item($r, 4).Value = ""
  $ws.Cells.Item($r, 5).Value = $inputPage
  $ws.Cells.Item($r, 6).Value = "not used"
  $ws.Cells.Item($r, 7).Value = "N/A"
  $ws.Cells.Item($r, 8).Value = "Ty"

  $r = $r + 1
}

$lastRow = $r - 1

# Widen column C slightly to fit the longer new names.
$ws.Columns.Item(3).ColumnWidth = 28

# Extend the "Type" dropdown validation down through the new rows.
$ws.Range("A2:A$lastRow").Validation.Delete()
$ws.Range("A2:A$lastRow").Validation.Add(3, 1, 1, "=Types")

# Update the view to match where editing left off.
$ws.Application.ActiveWindow.ScrollRow = 52
$ws.Range("C$lastRow").Select()
